$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sighting records in row 4 (Buskskvätta / Saxicola rubetra) and
# row 6 (Hällebräcka / Saxifraga osloënsis) had their data swapped.
# Only the columns that actually differ between the two rows need to be
# exchanged: Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Noggrannhet, plus the comment (AC) and the
# biotope description (AI), which move from row 6 to row 4.

$swapCols = @("A", "B", "D", "E", "F", "G", "H", "S")

foreach ($col in $swapCols) {
    $addr4 = $col + "4"
    $addr6 = $col + "6"

    $val4 = $ws.Range($addr4).Value2
    $val6 = $ws.Range($addr6).Value2

    $ws.Range($addr4).Value2 = $val6
    $ws.Range($addr6).Value2 = $val4
}

# "Publik kommentar" moves from row 6 to row 4.
$ws.Range("AC4").Value2 = $ws.Range("AC6").Value2
$ws.Range("AC6").ClearContents()

# "Biotop-beskrivning" moves from row 6 to row 4.
$ws.Range("AI4").Value2 = $ws.Range("AI6").Value2
$ws.Range("AI6").ClearContents()
